# Update Name of Algo
# Apply updated values to the KNN result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E7"   = 16.705
    "C9"   = -10.85
    "E12"  = 17.739
    "C18"  = -11.62
    "C20"  = -12.371
    "E26"  = 16.525
    "C27"  = -13.224
    "E27"  = 16.366
    "E29"  = 16.941
    "E37"  = 16.852
    "E38"  = 16.745
    "E51"  = 16.743
    "E55"  = 16.494
    "C69"  = -10.672
    "E69"  = 17.396
    "E70"  = 17.654
    "C76"  = -13.032
    "C82"  = -11.991
    "E83"  = 16.886
    "E102" = 16.784
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
